{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"856\u00f76=\", \"830\u00f79=\"],\n  [\"545\u00f76=\", \"367\u00f77=\"],\n  [\"697\u00f74=\", \"517\u00f75=\"],\n  [\"959\u00f79=\", \"293\u00f75=\"],\n  [\"808\u00f72=\", \"220\u00f78=\"],\n  [\"969\u00f74=\", \"806\u00f78=\"],\n  [\"951\u00f73=\", \"692\u00f75=\"],\n  [\"642\u00f74=\", \"646\u00f79=\"],\n  [\"263\u00f77=\", \"258\u00f72=\"],\n  [\"645\u00f79=\", \"321\u00f76=\"],\n  [\"697\u00f78=\", \"269\u00f79=\"],\n  [\"933\u00f76=\", \"244\u00f79=\"],\n  [\"783\u00f72=\", \"748\u00f77=\"],\n  [\"746\u00f79=\", \"785\u00f73=\"],\n  [\"577\u00f75=\", \"179\u00f72=\"],\n  [\"827\u00f79=\", \"505\u00f79=\"],\n  [\"231\u00f79=\", \"836\u00f79=\"],\n  [\"186\u00f77=\", \"129\u00f72=\"],\n  [\"825\u00f74=\", \"735\u00f76=\"],\n  [\"117\u00f77=\", \"459\u00f78=\"],\n  [\"193\u00f74=\", \"646\u00f78=\"],\n  [\"794\u00f74=\", \"517\u00f73=\"],\n  [\"214\u00f72=\", \"295\u00f79=\"],\n  [\"261\u00f77=\", \"934\u00f75=\"],\n  [\"411\u00f78=\", \"186\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"856\u00f76=\"\n$find.Replacement.Text = \"830\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"545\u00f76=\"\n$find.Replacement.Text = \"367\u00f77=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"697\u00f74=\"\n$find.Replacement.Text = \"517\u00f75=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"959\u00f79=\"\n$find.Replacement.Text = \"293\u00f75=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"808\u00f72=\"\n$find.Replacement.Text = \"220\u00f78=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"969\u00f74=\"\n$find.Replacement.Text = \"806\u00f78=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"951\u00f73=\"\n$find.Replacement.Text = \"692\u00f75=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"642\u00f74=\"\n$find.Replacement.Text = \"646\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"263\u00f77=\"\n$find.Replacement.Text = \"258\u00f72=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"645\u00f79=\"\n$find.Replacement.Text = \"321\u00f76=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"697\u00f78=\"\n$find.Replacement.Text = \"269\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"933\u00f76=\"\n$find.Replacement.Text = \"244\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"783\u00f72=\"\n$find.Replacement.Text = \"748\u00f77=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"746\u00f79=\"\n$find.Replacement.Text = \"785\u00f73=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"577\u00f75=\"\n$find.Replacement.Text = \"179\u00f72=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"827\u00f79=\"\n$find.Replacement.Text = \"505\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"231\u00f79=\"\n$find.Replacement.Text = \"836\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"186\u00f77=\"\n$find.Replacement.Text = \"129\u00f72=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"825\u00f74=\"\n$find.Replacement.Text = \"735\u00f76=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"117\u00f77=\"\n$find.Replacement.Text = \"459\u00f78=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"193\u00f74=\"\n$find.Replacement.Text = \"646\u00f78=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"794\u00f74=\"\n$find.Replacement.Text = \"517\u00f73=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"214\u00f72=\"\n$find.Replacement.Text = \"295\u00f79=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"261\u00f77=\"\n$find.Replacement.Text = \"934\u00f75=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n\n$find = $d.Content.Find\n$find.Text = \"411\u00f78=\"\n$find.Replacement.Text = \"186\u00f74=\"\n$find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n"}
